$wb = $excel.ActiveWorkbook

# Hunk 0: ALC!row9
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 300.75
$ws.Range("I9").Value = 300.75
$ws.Range("K9").Value = 300.75
$ws.Range("M9").Value = -131.75

# Hunk 1: ALC!row62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5140.6
$ws.Range("I62").Value = 5234.6665
$ws.Range("J62").Value = 4999.5
$ws.Range("K62").Value = 5234.6665
$ws.Range("L62").Value = 4999.5
$ws.Range("M62").Value = -4610.6665
$ws.Range("N62").Value = -6247.5

# Hunk 2: ALC!row65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 5140.6
$ws.Range("I65").Value = 5234.6665
$ws.Range("J65").Value = 4999.5
$ws.Range("K65").Value = 26173.3325
$ws.Range("L65").Value = 24997.5
$ws.Range("M65").Value = -23053.3325
$ws.Range("N65").Value = -31237.5

# Hunk 3: ALC!row100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5973.2856
$ws.Range("J100").Value = 8534
$ws.Range("L100").Value = 8534
$ws.Range("N100").Value = -9616

# Hunk 4: ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4329.6562
$ws.Range("I138").Value = 1683.1666
$ws.Range("J138").Value = 4940.385
$ws.Range("K138").Value = 5049.4998
$ws.Range("L138").Value = 14821.155
$ws.Range("M138").Value = 90.5002000000004
$ws.Range("N138").Value = -25101.155

# Hunk 5: ARM!row4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

# Hunk 6: ARM!row6
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()

# Hunk 7: ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5176.7104
$ws.Range("I32").Value = 5262.5835
$ws.Range("K32").Value = 5262.5835
$ws.Range("M32").Value = -4975.5835

# Hunk 8: ARM!row63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4480
$ws.Range("I63").Value = 4222
$ws.Range("J63").Value = 4867
$ws.Range("K63").Value = 4222
$ws.Range("L63").Value = 4867
$ws.Range("M63").Value = -3536
$ws.Range("N63").Value = -6239

# Hunk 9: ARM!row66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 4480
$ws.Range("I66").Value = 4222
$ws.Range("J66").Value = 4867
$ws.Range("K66").Value = 21110
$ws.Range("L66").Value = 24335
$ws.Range("M66").Value = -17678
$ws.Range("N66").Value = -31199

# Hunk 10: ARM!row92
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 51325
$ws.Range("J92").Value = 51325
$ws.Range("L92").Value = 51325
$ws.Range("N92").Value = -56317

# Hunk 11: ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1153.6562
$ws.Range("I132").Value = 1096.7333
$ws.Range("J132").Value = 2007.5
$ws.Range("K132").Value = 3290.199900000001
$ws.Range("L132").Value = 6022.5
$ws.Range("M132").Value = -760.1999000000005
$ws.Range("N132").Value = -11082.5

# Hunk 12: BSM!row88
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 34337
$ws.Range("J88").Value = 34337
$ws.Range("L88").Value = 34337
$ws.Range("N88").Value = -35149

# Hunk 13: BSM!row91
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H91").Value = 34337
$ws.Range("J91").Value = 34337
$ws.Range("L91").Value = 34337
$ws.Range("N91").Value = -37145

# Hunk 14: BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1831.9697
$ws.Range("I107").Value = 1900.2
$ws.Range("J107").Value = 1149.6666
$ws.Range("K107").Value = 1900.2
$ws.Range("L107").Value = 1149.6666
$ws.Range("M107").Value = 19.79999999999995
$ws.Range("N107").Value = -4989.6666

# Hunk 15: BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3735.125
$ws.Range("I134").Value = 3548.111
$ws.Range("K134").Value = 10644.333
$ws.Range("M134").Value = -8109.332999999999

# Hunk 16: CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2061.818
$ws.Range("I31").Value = 1919.7778
$ws.Range("J31").Value = 2701
$ws.Range("K31").Value = 1919.7778
$ws.Range("L31").Value = 2701
$ws.Range("M31").Value = -1624.7778
$ws.Range("N31").Value = -3291

# Hunk 17: CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2061.818
$ws.Range("I34").Value = 1919.7778
$ws.Range("J34").Value = 2701
$ws.Range("K34").Value = 1919.7778
$ws.Range("L34").Value = 2701
$ws.Range("M34").Value = -1717.7778
$ws.Range("N34").Value = -3105

# Hunk 18: CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 11815271
$ws.Range("I99").Value = 2443064
$ws.Range("J99").Value = 22228834
$ws.Range("K99").Value = 2443064
$ws.Range("L99").Value = 22228834
$ws.Range("M99").Value = -2441566
$ws.Range("N99").Value = -22231830

# Hunk 19: CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 11815271
$ws.Range("I126").Value = 2443064
$ws.Range("J126").Value = 22228834
$ws.Range("K126").Value = 7329192
$ws.Range("L126").Value = 66686502
$ws.Range("M126").Value = -7326722
$ws.Range("N126").Value = -66691442

# Hunk 20: CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3789.9119
$ws.Range("I134").Value = 3002.5217
$ws.Range("K134").Value = 9007.5651
$ws.Range("M134").Value = -6472.5651

# Hunk 21: GSM!row43
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2583.2
$ws.Range("I43").Value = 2583.2
$ws.Range("K43").Value = 2583.2
$ws.Range("M43").Value = -2432.2

# Hunk 22: GSM!row44
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 21983.334
$ws.Range("I44").Value = 21975
$ws.Range("J44").Value = 21987.5
$ws.Range("K44").Value = 21975
$ws.Range("L44").Value = 21987.5
$ws.Range("M44").Value = -21379
$ws.Range("N44").Value = -23179.5

# Hunk 23: GSM!row47
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 22499.5
$ws.Range("J47").Value = 22499.5
$ws.Range("L47").Value = 22499.5
$ws.Range("N47").Value = -23635.5

# Hunk 24: GSM!row52
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 40250
$ws.Range("J52").Value = 40250
$ws.Range("L52").Value = 40250
$ws.Range("N52").Value = -40768

# Hunk 25: GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 83230.13
$ws.Range("I70").Value = 192169.5
$ws.Range("J70").Value = 10603.889
$ws.Range("K70").Value = 192169.5
$ws.Range("L70").Value = 10603.889
$ws.Range("M70").Value = -191899.5
$ws.Range("N70").Value = -11143.889

# Hunk 26: GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 83230.13
$ws.Range("I73").Value = 192169.5
$ws.Range("J73").Value = 10603.889
$ws.Range("K73").Value = 192169.5
$ws.Range("L73").Value = 10603.889
$ws.Range("M73").Value = -191233.5
$ws.Range("N73").Value = -12475.889

# Hunk 27: GSM!row107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 738.6667
$ws.Range("I107").Value = 958.75
$ws.Range("J107").Value = 298.5
$ws.Range("K107").Value = 958.75
$ws.Range("L107").Value = 298.5
$ws.Range("M107").Value = 961.25
$ws.Range("N107").Value = -4138.5

# Hunk 28: GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1944.1333
$ws.Range("I132").Value = 1906.6923
$ws.Range("J132").Value = 2187.5
$ws.Range("K132").Value = 5720.0769
$ws.Range("L132").Value = 6562.5
$ws.Range("M132").Value = -3190.0769
$ws.Range("N132").Value = -11622.5

# Hunk 29: LTW!row60
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 27000
$ws.Range("I60").Value = 14000
$ws.Range("J60").Value = 40000
$ws.Range("K60").Value = 14000
$ws.Range("L60").Value = 40000
$ws.Range("M60").Value = -13491
$ws.Range("N60").Value = -41018

# Hunk 30: LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4996.875
$ws.Range("I132").Value = 3514.1667
$ws.Range("K132").Value = 10542.5001
$ws.Range("M132").Value = -8012.500100000001

# Hunk 31: LTW!row139
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H139").Value = 98750
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()

# Hunk 32: WVR!row14
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 8876.125
$ws.Range("I14").Value = 8202
$ws.Range("J14").Value = 9999.666999999999
$ws.Range("K14").Value = 8202
$ws.Range("L14").Value = 9999.666999999999
$ws.Range("M14").Value = -8034
$ws.Range("N14").Value = -10335.667

# Hunk 33: WVR!row58
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 33747.5
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 33747.5
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 33747.5
$ws.Range("N58").Value = -34363.5
$ws.Range("M58").ClearContents()

# Hunk 34: WVR!row103
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Hunk 35: WVR!row113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 530.3333
$ws.Range("I113").Value = 431.65216
$ws.Range("K113").Value = 1294.95648
$ws.Range("M113").Value = 875.0435200000002
